$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.844.02"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.308.91"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'538.53"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").Value = "'132.37"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("D9").Value = "2.306.98"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "'0.101"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "'5.49"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "'23.82"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "2.718.63"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "58.697.24"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "2.297.94"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'10.64"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'4.17"
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("D21").Value = "'313.13"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "'6.63"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'62.37"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'7.94"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.72"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'170.70"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "0.0₃0734"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'17.91"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +4.68%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'4.07"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("D40").Value = "'1.52"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "'297.97"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").Value = "'141.60"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'3.45"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "'0.0960"
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").Value = "'0.0496"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'18.31"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +1.03%  "
